$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# Append the new log row (row 22)
$row = 22
$ws.Cells.Item($row, 1).Value = "Nieuwe bestelling"
$ws.Cells.Item($row, 2).Value = "planning@testbedrijf123.nl"
$ws.Cells.Item($row, 3).Value = "Wil je 200 stuks M8-bouten bestellen bij onze leverancier?"
$ws.Cells.Item($row, 4).Value = "Intern verzoek / Actie voor medewerker"
$ws.Cells.Item($row, 5).Value = "Bedankt, we hebben dit doorgestuurd naar inkoop@testbedrijf123.nl."
$ws.Cells.Item($row, 6).Value = "2025-08-14 21:15:54"
$ws.Cells.Item($row, 7).Value = "Nee"
$ws.Cells.Item($row, 8).Value = "Ja"
$ws.Cells.Item($row, 9).Value = "Nee"
$ws.Cells.Item($row, 10).Value = "Nee"

# Extend the conditional-formatting ranges to cover the new row
$cols = @("D", "G", "H", "I", "J")
foreach ($col in $cols) {
    $oldRange = $ws.Range($col + "2:" + $col + "21")
    $newRange = $ws.Range($col + "2:" + $col + "22")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# Update the dashboard summary count for the category that gained a row
$dash.Cells.Item(2, 2).Value = 16
